$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(1201, 2, 10, 10, 10, 10)
  3 = @(401, 9, 48, 67, 75, 45)
  4 = @(801, 3, 67, 65, 52, 45)
  5 = @(601, 9, 60, 67, 60, 42)
  6 = @(1203, 3, 15, 15, 15, 15)
  7 = @(902, 1, 0, 0, 0, 0)
  8 = @(1001, 18, 30, 75, 60, 72)
  9 = @(301, 6, 45, 30, 60, 45)
  10 = @(501, 9, 52, 30, 75, 45)
  11 = @(701, 3, 90, 45, 97, 15)
  12 = @(1202, 2, 10, 10, 10, 10)
  13 = @(101, 9, 30, 15, 60, 15)
  14 = @(901, 16, 15, 45, 60, 60)
  15 = @(201, 9, 30, 15, 45, 30)
  16 = @(2, 0, 2, 2, 2, 2)
  17 = @(502, 0, 4, 0, 0, 0)
  18 = @(1101, 0, 15, 30, 30, 0)
  19 = @(3, 0, 3, 3, 3, 3)
  20 = @(802, 0, 4, 5, 4, 0)
  21 = @(1, 0, 2, 2, 2, 2)
  22 = @(602, 0, 0, 4, 0, 9)
  23 = @(402, 0, 0, 4, 0, 0)
}

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  for ($i = 0; $i -lt $rowVals.Length; $i++) {
    $ws.Cells.Item([int]$r, $i + 1).Value = $rowVals[$i]
  }
}
